# #327 Ajout des profils d'acces
# - Bump the StructureDefinition "Date" metadata value.
# - Swap the two "Mapping" columns (AK/AL) on the Elements sheet:
#   "RIM Mapping" and "Spécification métier vers l'extension ROR
#   LocationResidentialCapacity" traded places, header included.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the generation date ---------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- Elements sheet: swap the AK / AL "Mapping" columns --------------
$wsElem = $wb.Worksheets.Item("Elements")

# Column headers (row 1)
$wsElem.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR LocationResidentialCapacity"
$wsElem.Range("AL1").Value = "Mapping: RIM Mapping"

# Data rows: swap AK<->AL cell by cell (rows whose AK/AL are both
# already blank are left untouched on purpose).
$rows = 2..16
foreach ($r in $rows) {
    $akCell = $wsElem.Range("AK$r")
    $alCell = $wsElem.Range("AL$r")
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Column widths followed the content: AK becomes the wide
# "Spécification métier" column, AL becomes the narrow "RIM Mapping" one.
$wsElem.Columns.Item(37).ColumnWidth = 83.55338541666667
$wsElem.Columns.Item(38).ColumnWidth = 24.147135416666668
